$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 5.652167666666667
$ws.Cells.Item(2,8).Value2 = 16.956503
$ws.Cells.Item(2,9).Value2 = 0.1860329065948871
$ws.Cells.Item(2,10).Value2 = 0.1860329065948871
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 3.211751333333333
$ws.Cells.Item(2,14).Value2 = 9.635254
$ws.Cells.Item(2,15).Value2 = 0.9274105493513785
$ws.Cells.Item(2,16).Value2 = 0.9274105493513782
$ws.Cells.Item(2,17).Value2 = 18.15335703964022
$ws.Cells.Item(2,18).Value2 = 163.380213356762
$ws.Cells.Item(2,19).Value2 = 0.1725288801025979
$ws.Cells.Item(2,20).Value2 = 0.1725288801025978

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 5.652167666666667
$ws.Cells.Item(3,8).Value2 = 16.956503
$ws.Cells.Item(3,9).Value2 = 0.1860329065948871
$ws.Cells.Item(3,10).Value2 = 0.1860329065948871
$ws.Cells.Item(3,11).Value2 = 2
$ws.Cells.Item(3,12).Value2 = 0.6666666666666666
$ws.Cells.Item(3,13).Value2 = 0.2513873333333334
$ws.Cells.Item(3,14).Value2 = 0.754162
$ws.Cells.Item(3,15).Value2 = 0.07258945064862164
$ws.Cells.Item(3,16).Value2 = 0.07258945064862163
$ws.Cells.Item(3,17).Value2 = 1.420883357276223
$ws.Cells.Item(3,18).Value2 = 12.787950215486
$ws.Cells.Item(3,19).Value2 = 0.01350402649228919
$ws.Cells.Item(3,20).Value2 = 0.01350402649228919

# Row 4: FAPs -> ECs
$ws.Cells.Item(4,1).Value2 = "FAPs"
$ws.Cells.Item(4,4).Value2 = "ECs"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 16.407289
$ws.Cells.Item(4,8).Value2 = 49.221867
$ws.Cells.Item(4,9).Value2 = 0.5400221369958743
$ws.Cells.Item(4,10).Value2 = 0.5400221369958743
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 3.211751333333333
$ws.Cells.Item(4,14).Value2 = 9.635254
$ws.Cells.Item(4,15).Value2 = 0.9274105493513785
$ws.Cells.Item(4,16).Value2 = 0.9274105493513782
$ws.Cells.Item(4,17).Value2 = 52.69613232213534
$ws.Cells.Item(4,18).Value2 = 474.265190899218
$ws.Cells.Item(4,19).Value2 = 0.5008222267332492
$ws.Cells.Item(4,20).Value2 = 0.5008222267332491

# Row 5 (new): FAPs -> FAPs
$ws.Cells.Item(5,1).Value2 = "FAPs"
$ws.Cells.Item(5,2).Value2 = "Lrpap1"
$ws.Cells.Item(5,3).Value2 = "Lrp8"
$ws.Cells.Item(5,4).Value2 = "FAPs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 16.407289
$ws.Cells.Item(5,8).Value2 = 49.221867
$ws.Cells.Item(5,9).Value2 = 0.5400221369958743
$ws.Cells.Item(5,10).Value2 = 0.5400221369958743
$ws.Cells.Item(5,11).Value2 = 2
$ws.Cells.Item(5,12).Value2 = 0.6666666666666666
$ws.Cells.Item(5,13).Value2 = 0.2513873333333334
$ws.Cells.Item(5,14).Value2 = 0.754162
$ws.Cells.Item(5,15).Value2 = 0.07258945064862164
$ws.Cells.Item(5,16).Value2 = 0.07258945064862163
$ws.Cells.Item(5,17).Value2 = 4.124584628939334
$ws.Cells.Item(5,18).Value2 = 37.121261660454
$ws.Cells.Item(5,19).Value2 = 0.03919991026262521
$ws.Cells.Item(5,20).Value2 = 0.0391999102626252

# Row 6 (new): sCs -> ECs
$ws.Cells.Item(6,1).Value2 = "sCs"
$ws.Cells.Item(6,2).Value2 = "Lrpap1"
$ws.Cells.Item(6,3).Value2 = "Lrp8"
$ws.Cells.Item(6,4).Value2 = "ECs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 8.323166333333333
$ws.Cells.Item(6,8).Value2 = 24.969499
$ws.Cells.Item(6,9).Value2 = 0.2739449564092387
$ws.Cells.Item(6,10).Value2 = 0.2739449564092387
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 3.211751333333333
$ws.Cells.Item(6,14).Value2 = 9.635254
$ws.Cells.Item(6,15).Value2 = 0.9274105493513785
$ws.Cells.Item(6,16).Value2 = 0.9274105493513782
$ws.Cells.Item(6,17).Value2 = 26.73194056863844
$ws.Cells.Item(6,18).Value2 = 240.587465117746
$ws.Cells.Item(6,19).Value2 = 0.2540594425155315
$ws.Cells.Item(6,20).Value2 = 0.2540594425155314

# Row 7 (new): sCs -> FAPs
$ws.Cells.Item(7,1).Value2 = "sCs"
$ws.Cells.Item(7,2).Value2 = "Lrpap1"
$ws.Cells.Item(7,3).Value2 = "Lrp8"
$ws.Cells.Item(7,4).Value2 = "FAPs"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 8.323166333333333
$ws.Cells.Item(7,8).Value2 = 24.969499
$ws.Cells.Item(7,9).Value2 = 0.2739449564092387
$ws.Cells.Item(7,10).Value2 = 0.2739449564092387
$ws.Cells.Item(7,11).Value2 = 2
$ws.Cells.Item(7,12).Value2 = 0.6666666666666666
$ws.Cells.Item(7,13).Value2 = 0.2513873333333334
$ws.Cells.Item(7,14).Value2 = 0.754162
$ws.Cells.Item(7,15).Value2 = 0.07258945064862164
$ws.Cells.Item(7,16).Value2 = 0.07258945064862163
$ws.Cells.Item(7,17).Value2 = 2.092338589426444
$ws.Cells.Item(7,18).Value2 = 18.831047304838
$ws.Cells.Item(7,19).Value2 = 0.01988551389370724
$ws.Cells.Item(7,20).Value2 = 0.01988551389370724
